$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Target values for B2:F261 (one comma-separated row per spreadsheet row;
# an empty field means that cell has no value in the workbook and must be
# left blank, not zero).
$csvData = @"
3,3,3,2,3
9,12,8,7,12
12,5,4,8,5
8,7,6,6,9
2,6,9,2,2
10,1,8,5,6
9,12,7,7,12
6,8,4,9,9
4,10,2,8,10
1,2,8,3,1
7,4,12,9,12
9,8,4,12,6
10,6,8,8,10
12,10,6,2,3
1,1,4,12,5
10,6,12,4,12
9,12,3,11,10
5,7,8,6,6
4,9,9,3,3
2,0,10,12,8
6,6,4,6,5
5,8,7,7,10
12,12,6,1,12
3,1,8,11,3
10,5,1,10,8
4,6,5,7,12
6,12,12,1,9
12,1,4,10,7
2,10,6,4,2
5,7,12,7,9
10,5,9,0,12
9,1,10,12,10
1,7,3,9,5
12,8,6,7,0
9,12,4,5,9
10,9,9,0,8
,3,10,10,12
,12,8,7,4
,7,6,12,0
,8,10,3,6
,3,0,10,4
,12,8,12,8
,4,4,7,12
,10,10,3,2
,1,9,4,6
,8,11,12,9
,6,8,7,4
,7,5,2,8
,3,3,4,3
,8,4,8,12
,4,9,6,9
,7,5,3,8
,5,6,9,1
,0,9,8,5
,8,4,7,9
,10,5,3,4
,12,2,4,3
,2,4,6,6
,8,9,12,5
,12,5,3,8
,4,8,5,12
,0,10,7,2
,9,5,8,6
,6,7,12,9
,4,2,1,7
,7,10,5,5
,2,11,10,2
,6,8,6,10
,8,6,3,6
,12,9,5,12
,3,10,9,3
,4,12,4,5
,5,1,2,4
,8,9,10,12
,3,12,8,3
,6,11,4,9
,4,6,11,4
,5,9,1,6
,6,11,10,8
,0,6,5,3
,12,0,12,4
,10,8,9,12
,8,5,0,9
,1,7,12,2
,11,8,7,6
,4,11,11,8
,8,4,1,9
,1,7,7,5
,5,3,4,3
,8,4,5,8
,12,7,1,12
,3,6,8,5
,9,8,5,4
,11,12,9,1
,4,10,3,6
,3,11,4,9
,6,1,11,8
,8,10,8,12
,10,6,7,1
,2,7,1,10
,4,11,9,12
,6,12,4,5
,8,10,11,10
,3,6,3,2
,7,0,9,8
,8,8,11,12
,6,12,7,5
,3,4,2,3
,5,7,5,10
,12,6,12,5
,4,8,7,7
,2,12,4,0
,12,3,2,9
,8,6,6,12
,6,12,12,4
,5,10,4,6
,3,11,5,2
,12,7,2,8
,4,10,11,4
,8,12,8,6
,2,3,4,0
,10,4,1,7
,5,12,5,6
,4,5,12,9
,8,4,10,0
,1,8,3,4
,5,5,6,6
,12,9,10,12
,8,3,9,1
,2,6,11,8
,5,4,1,4
,6,10,4,6
,12,7,11,3
,11,4,10,7
,2,6,2,12
,5,12,5,9
,4,0,7,3
,7,4,4,5
,5,7,1,9
,1,8,5,12
,10,10,4,4
,12,11,7,2
,6,6,2,6
,11,12,8,12
,3,1,6,4
,5,10,9,8
,6,9,12,3
,8,12,1,4
,2,6,5,6
,12,8,12,10
,9,7,6,12
,4,12,0,3
,2,2,4,9
,12,4,5,5
,10,10,9,6
,6,7,1,2
,5,6,12,9
,3,8,7,12
,11,12,8,6
,12,10,3,9
,6,1,5,
,5,7,12,
,3,4,7,
,4,8,3,
,12,6,6,
,8,12,10,
,5,8,5,
,3,6,8,
,8,2,2,
,6,12,9,
,12,6,10,
,0,10,8,
,5,12,3,
,12,8,6,
,8,4,5,
,5,10,7,
,,3,10,
,,12,2,
,,10,7,
,,7,9,
,,9,12,
,,8,10,
,,10,2,
,,5,7,
,,3,6,
,,10,4,
,,12,1,
,,4,8,
,,10,6,
,,12,7,
,,4,1,
,,8,9,
,,0,4,
,,4,8,
,,9,3,
,,6,12,
,,11,10,
,,12,5,
,,6,1,
,,10,8,
,,,7,
,,,5,
,,,3,
,,,12,
,,,7,
,,,8,
,,,3,
,,,5,
,,,7,
,,,6,
,,,3,
,,,12,
,,,9,
,,,7,
,,,2,
,,,4,
,,,10,
,,,7,
,,,6,
,,,3,
,,,8,
,,,4,
,,,6,
,,,7,
,,,1,
,,,5,
,,,4,
,,,12,
,,,9,
,,,2,
,,,6,
,,,9,
,,,5,
,,,6,
,,,3,
,,,7,
,,,10,
,,,4,
,,,0,
,,,9,
,,,12,
,,,7,
,,,1,
,,,12,
,,,7,
,,,4,
,,,1,
,,,10,
,,,9,
,,,5,
,,,1,
,,,7,
,,,8,
,,,6,
,,,10,
,,,1,
,,,5,
,,,10,
,,,4,
,,,7,
"@

$rowsText = $csvData -split "`n"
$numRows = $rowsText.Count
$numCols = 5

$arr = New-Object 'object[,]' $numRows,$numCols
for ($i = 0; $i -lt $numRows; $i++) {
    $fields = $rowsText[$i] -split ","
    for ($j = 0; $j -lt $numCols; $j++) {
        $arr[$i,$j] = $fields[$j]
    }
}

$ws.Range("B2:F261").Value = $arr
Write-Host "Updated $numRows rows"
